# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing team-statistics table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting already used by the other header cells
# (bold font + thin border + centered/top alignment) onto the three new
# header cells, then set their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is constant for every player row on this sheet (2-48).
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 74  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
